$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 700
$ws.Range("B3").Value = 300
$ws.Range("B5").Value = 250
$ws.Range("B6").Value = 150
$ws.Range("B7").Value = 350
$ws.Range("B8").Value = 150
